$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the new longer row labels
$ws.Columns.Item(1).ColumnWidth = 25.166666666666668

# Row 24: 8 yilda max hektar (max across the whole 8-year data block)
$ws.Range("A24").Value = "8 yılda max hektar"
$ws.Range("B24").Formula = "=MAX(B2:I16)"
$ws.Range("B24").NumberFormat = "#,##0"

# Row 25: 8 yilda min hektar (min across the whole 8-year data block)
$ws.Range("A25").Value = "8 yılda min hektar"
$ws.Range("B25").Formula = "=MIN(B2:I16)"
$ws.Range("B25").NumberFormat = "#,##0"

# Row 23: 2015 toplam hektar  (sum of 2015 column I, rows 2-16)
$ws.Range("A23").Value = "2015 toplam hektar "
$ws.Range("B23").Formula = "=SUM(I2:I16)"

# Row 26: 2015 ortalama hektar (average of 2015 column I, rows 2-16)
$ws.Range("A26").Value = "2015 ortalama hektar"
$ws.Range("B26").Formula = "=AVERAGE(I2:I16)"

# Row 27: 8 yilda ortalama hektar (average of the TOPLAM row across years)
$ws.Range("A27").Value = "8 yılda ortalama hektar"
$ws.Range("B27").Formula = "=AVERAGE(B17:I17)"

# Leftover formatting touches from the author (no content, just number formats)
$ws.Range("B28").NumberFormat = "0.000000%"
$ws.Range("D23").NumberFormat = "0.0000000"

# Final selection left on A28 after the edits
$ws.Range("A28").Select() | Out-Null
